# AP z boson CMS 8 TeV
# Add a new "process" column (Y) describing the physics process for every
# data row on the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell
$ws.Range("Y1").Value = "process"

# New data cells (rows 2-35) all share the same process description
$ws.Range("Y2:Y35").Value = "pp->Z/gamma*->l+ l-"

# Reflect the new column in the view: select it, just like after typing
# the values in the real editing session
[void]$ws.Range("Y1:Y35").Select()
